$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Add the two new worksheets at the end of the workbook, preserving order:
#   ... StatOutput, StatOutput_Message, CaseDetailStat, CaseDetailStat_Message
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$wsStat = $wb.Worksheets.Add($null, $lastSheet)
$wsStat.Name = "CaseDetailStat"

$wsMsg = $wb.Worksheets.Add($null, $wsStat)
$wsMsg.Name = "CaseDetailStat_Message"

# ---------------------------------------------------------------------------
# CaseDetailStat sheet: header row + one data row describing the file found
# for case NCATS-COP01CCB010072.
# ---------------------------------------------------------------------------
$wsStat.Range("A1").Value = "File Name"
$wsStat.Range("B1").Value = "File Type"
$wsStat.Range("C1").Value = "Association"
$wsStat.Range("D1").Value = "Description"
$wsStat.Range("E1").Value = "Format"
$wsStat.Range("F1").Value = "Size"

$wsStat.Range("A2").Value = "CCB010072.pdf"
$wsStat.Range("B2").Value = "Pathology Report"
$wsStat.Range("C2").Value = "diagnosis"
$wsStat.Range("D2").Value = ""
$wsStat.Range("E2").Value = "pdf"
# Leading apostrophe forces this numeric-looking value to be stored as text
# (matches the source data, which is a text cell, not a real number).
$wsStat.Range("F2").Value = "'57.732421875"

# ---------------------------------------------------------------------------
# CaseDetailStat_Message sheet: three repeated "log" blocks of
#   Neo4j_URL: / <url>
#   User_name: / neo4j
#   PWD: / <pwd>
#   Cypher: / <cypher text>
#   Output: / <output file path>
# The Cypher text differs per block; the rest repeats verbatim.
# ---------------------------------------------------------------------------
$neo4jUrl = "bolt://ncias-q2251-c.nci.nih.gov:7687"
$userName = "neo4j"
$pwd_ = "icdcDBneo4j0"
$outputPath = 'C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC04_Canine_Filter_Breed-BassHnd_Neo4jData.xlsx'

$cypher1 = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN [''Basset Hound''] WITH DISTINCT c AS c, p, s, demo, diag RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(s.clinical_study_designation,'''') AS `Study Code` , coalesce(s.clinical_study_type,'''') AS  `Study Type`, coalesce(demo.breed,'''') AS Breed , coalesce(diag.disease_term,'''') AS Diagnosis , coalesce(diag.stage_of_disease,'''') AS `Stage of Disease` ,  coalesce(demo.patient_age_at_enrollment,'''') AS Age , coalesce(demo.sex,'''') AS Sex , coalesce(demo.neutered_indicator,'''') AS  `Neutered Status`'

$cypher2 = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN[''Basset Hound'']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study'

$cypher3 = 'MATCH (f:file)-[*]->(c:case) WITH DISTINCT(f) AS f, c MATCH (f)-->(parent) WHERE c.case_id IN [''NCATS-COP01CCB010072''] RETURN f.file_name AS `File Name` ,f.file_type AS `File Type`,head(labels(parent)) AS `Association`, f.file_description AS `Description`,f.file_format AS Format,((f.file_size)/1024) AS Size'

$rows = @(
  "Neo4j_URL:", $neo4jUrl, "User_name:", $userName, "PWD:", $pwd_, "Cypher:", $cypher1, "Output:", $outputPath,
  "Neo4j_URL:", $neo4jUrl, "User_name:", $userName, "PWD:", $pwd_, "Cypher:", $cypher2, "Output:", $outputPath,
  "Neo4j_URL:", $neo4jUrl, "User_name:", $userName, "PWD:", $pwd_, "Cypher:", $cypher3, "Output:", $outputPath
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 1
    $wsMsg.Cells.Item($r, 1).Value = $rows[$i]
}

Write-Output "Added CaseDetailStat and CaseDetailStat_Message sheets"
